# Applies the "Updated cryptos list" data refresh described in the commit diff.
# All D/E (and the B/C swap for rows 39-40) cells are stored as *text* in the
# sheet (inline strings in the source OOXML), even though many look numeric
# (e.g. "600.92", "0.140"). Plain `Range.Value = "600.92"` would let Excel
# auto-coerce that into a real number, so each write forces the cell to Text
# ("@") first, and then resets the cell style back to "Normal" afterwards so
# no stray number-format/style is left behind on the cell (matches the source
# diff, which touches only cell text, not styles).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Address, $Value)
    $r = $ws.Range($Address)
    $r.NumberFormat = "@"
    $r.Value = $Value
    $r.Style = "Normal"
}

Set-TextValue "D2" "68.318.81"
Set-TextValue "E2" "  -2.26%  "
Set-TextValue "D3" "3.823.25"
Set-TextValue "E3" "  -2.36%  "
Set-TextValue "E4" "  +0.00%  "
Set-TextValue "D5" "600.92"
Set-TextValue "E5" "  -1.21%  "
Set-TextValue "D6" "169.52"
Set-TextValue "E6" "  -0.26%  "
Set-TextValue "D7" "3.830.60"
Set-TextValue "E7" "  -2.15%  "
Set-TextValue "E8" "  -0.19%  "
Set-TextValue "D9" "0.527"
Set-TextValue "E9" "  -1.80%  "
Set-TextValue "E10" "  -2.97%  "
Set-TextValue "E11" "  +0.96%  "
Set-TextValue "D12" "0.458"
Set-TextValue "E12" "  -2.53%  "
Set-TextValue "D13" "0.0000264"
Set-TextValue "E13" "  +3.01%  "
Set-TextValue "D14" "37.07"
Set-TextValue "E14" "  -3.43%  "
Set-TextValue "D15" "4.461.38"
Set-TextValue "E15" "  -2.50%  "
Set-TextValue "D16" "3.823.31"
Set-TextValue "E16" "  -2.74%  "
Set-TextValue "D17" "68.279.59"
Set-TextValue "E17" "  -2.37%  "
Set-TextValue "D18" "18.46"
Set-TextValue "E18" "  -1.91%  "
Set-TextValue "D19" "7.41"
Set-TextValue "E19" "  -2.95%  "
Set-TextValue "E20" "  -0.78%  "
Set-TextValue "D21" "11.10"
Set-TextValue "E21" "  -1.07%  "
Set-TextValue "D22" "470.60"
Set-TextValue "E22" "  -4.62%  "
Set-TextValue "D23" "0.734"
Set-TextValue "E23" "  -1.95%  "
Set-TextValue "D24" "0.0000161"
Set-TextValue "E24" "  -4.54%  "
Set-TextValue "D25" "83.13"
Set-TextValue "E25" "  -2.95%  "
Set-TextValue "E26" "  -3.28%  "
Set-TextValue "D27" "12.16"
Set-TextValue "E27" "  -1.90%  "
Set-TextValue "E28" "  -1.72%  "
Set-TextValue "E29" "  +0.03%  "
Set-TextValue "D30" "2.96"
Set-TextValue "E30" "  -1.18%  "
Set-TextValue "D31" "3.969.06"
Set-TextValue "E31" "  -2.45%  "
Set-TextValue "D32" "7.68"
Set-TextValue "E32" "  -2.31%  "
Set-TextValue "D33" "31.55"
Set-TextValue "E33" "  -1.89%  "
Set-TextValue "E34" "  -5.29%  "
Set-TextValue "D35" "9.45"
Set-TextValue "E35" "  -1.57%  "
Set-TextValue "D36" "3.782.89"
Set-TextValue "E36" "  -2.54%  "
Set-TextValue "E37" "  -2.97%  "
Set-TextValue "E38" "  +10.75%  "
Set-TextValue "B39" "Mantle"
Set-TextValue "C39" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D39" "1.02"
Set-TextValue "E39" "  -2.69%  "
Set-TextValue "B40" "Kaspa"
Set-TextValue "C40" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D40" "0.140"
Set-TextValue "E40" "  -1.59%  "
Set-TextValue "D41" "5.94"
Set-TextValue "E41" "  -3.31%  "
Set-TextValue "E42" "  -0.01%  "
Set-TextValue "D43" "0.315"
Set-TextValue "E43" "  -4.77%  "
Set-TextValue "E44" "  -7.03%  "
Set-TextValue "E45" "  +0.56%  "
Set-TextValue "D46" "0.000296"
Set-TextValue "E46" "  +9.01%  "
Set-TextValue "D48" "415.92"
Set-TextValue "E48" "  -5.44%  "
Set-TextValue "D49" "47.03"
Set-TextValue "E49" "  -2.46%  "
Set-TextValue "E50" "  +3.65%  "
Set-TextValue "D51" "141.28"
Set-TextValue "E51" "  -1.86%  "
